$wb = $excel.ActiveWorkbook

# 1. Update the status text from "Ready for handoff" to "In Translation"
#    wherever it appears (Overview!E2:F2, zh-cn!C2, de-de!C2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# 2. Narrow the "Status" columns to match the new (shorter) content.
#    The Excel column-width model stores a "pixel padded" character width
#    (stored width = ColumnWidth + ~0.8333 char, snapped to pixel grid), so
#    to land the persisted <col width="..."/> as close as possible to the
#    target 13.4101845877511 we need to set ColumnWidth to the target minus
#    that padding constant (~12.5 lands in the correct pixel bucket).
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
